$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 17
$ws.Range("E2").Value = 17.6299991607666
$ws.Range("F2").Value = 20.60000038146973
$ws.Range("G2").Value = 15.5
$ws.Range("H2").Value = 328636590
$ws.Range("I2").Value = "PSTG"

$ws.Range("D3").Value = 17
$ws.Range("E3").Value = 17.6299991607666
$ws.Range("F3").Value = 20.60000038146973
$ws.Range("G3").Value = 15.5
$ws.Range("H3").Value = 328636590
$ws.Range("I3").Value = "PSTG"

$ws.Range("D4").Value = 17
$ws.Range("E4").Value = 17.6299991607666
$ws.Range("F4").Value = 20.60000038146973
$ws.Range("G4").Value = 15.5
$ws.Range("H4").Value = 328636590
$ws.Range("I4").Value = "PSTG"

$ws.Range("D5").Value = 18.35000038146973
$ws.Range("E5").Value = 12.90999984741211
$ws.Range("F5").Value = 18.35899925231934
$ws.Range("G5").Value = 12.26000022888184
$ws.Range("H5").Value = 328636590
$ws.Range("I5").Value = "PSTG"

$ws.Range("D6").Value = 12.96000003814697
$ws.Range("E6").Value = 14.38000011444092
$ws.Range("F6").Value = 14.47999954223633
$ws.Range("G6").Value = 11.05000019073486
$ws.Range("H6").Value = 328636590
$ws.Range("I6").Value = "PSTG"

$ws.Range("D7").Value = 14.51000022888184
$ws.Range("E7").Value = 11.77999973297119
$ws.Range("F7").Value = 15.19999980926514
$ws.Range("G7").Value = 11.05000019073486
$ws.Range("H7").Value = 328636590
$ws.Range("I7").Value = "PSTG"

$ws.Range("D8").Value = 12.60000038146973
$ws.Range("E8").Value = 11.71000003814697
$ws.Range("F8").Value = 13.97000026702881
$ws.Range("G8").Value = 11
$ws.Range("H8").Value = 328636590
$ws.Range("I8").Value = "PSTG"

$ws.Range("D9").Value = 12.32999992370606
$ws.Range("E9").Value = 13.94999980926514
$ws.Range("F9").Value = 15.14000034332275
$ws.Range("G9").Value = 11.4399995803833
$ws.Range("H9").Value = 328636590
$ws.Range("I9").Value = "PSTG"

$ws.Range("D10").Value = 11.36999988555908
$ws.Range("E10").Value = 11.39999961853027
$ws.Range("F10").Value = 12.34000015258789
$ws.Range("G10").Value = 11.0600004196167
$ws.Range("H10").Value = 328636590
$ws.Range("I10").Value = "PSTG"

$ws.Range("D11").Value = 10.60000038146973
$ws.Range("E11").Value = 12.9399995803833
$ws.Range("F11").Value = 13.10000038146973
$ws.Range("G11").Value = 9.810000419616699
$ws.Range("H11").Value = 328636590
$ws.Range("I11").Value = "PSTG"

$ws.Range("D12").Value = 12.11999988555908
$ws.Range("E12").Value = 14.89000034332275
$ws.Range("F12").Value = 15.09000015258789
$ws.Range("G12").Value = 12
$ws.Range("H12").Value = 328636590
$ws.Range("I12").Value = "PSTG"

$ws.Range("D13").Value = 16.6299991607666
$ws.Range("E13").Value = 18.47999954223633
$ws.Range("F13").Value = 19.27499961853028
$ws.Range("G13").Value = 16.28000068664551
$ws.Range("H13").Value = 328636590
$ws.Range("I13").Value = "PSTG"

$ws.Range("D14").Value = 20.10000038146973
$ws.Range("E14").Value = 21.67000007629395
$ws.Range("F14").Value = 22.60000038146973
$ws.Range("G14").Value = 18.3700008392334
$ws.Range("H14").Value = 328636590
$ws.Range("I14").Value = "PSTG"

$ws.Range("D15").Value = 20.13999938964844
$ws.Range("E15").Value = 21.45999908447266
$ws.Range("F15").Value = 24.30999946594238
$ws.Range("G15").Value = 19.6200008392334
$ws.Range("H15").Value = 328636590
$ws.Range("I15").Value = "PSTG"

$ws.Range("D16").Value = 21.77000045776367
$ws.Range("E16").Value = 26.84000015258789
$ws.Range("F16").Value = 27.14999961853028
$ws.Range("G16").Value = 21.07999992370605
$ws.Range("H16").Value = 328636590
$ws.Range("I16").Value = "PSTG"

$ws.Range("D17").Value = 20.27000045776367
$ws.Range("E17").Value = 18.90999984741211
$ws.Range("F17").Value = 22.5
$ws.Range("G17").Value = 16.56999969482422
$ws.Range("H17").Value = 328636590
$ws.Range("I17").Value = "PSTG"

$ws.Range("D18").Value = 17.85000038146973
$ws.Range("E18").Value = 20.47999954223633
$ws.Range("F18").Value = 20.75
$ws.Range("G18").Value = 17.85000038146973
$ws.Range("H18").Value = 328636590
$ws.Range("I18").Value = "PSTG"

$ws.Range("D19").Value = 23.19000053405762
$ws.Range("E19").Value = 15.85999965667725
$ws.Range("F19").Value = 23.29999923706055
$ws.Range("G19").Value = 15.17099952697754
$ws.Range("H19").Value = 328636590
$ws.Range("I19").Value = "PSTG"

$ws.Range("D20").Value = 15.21000003814697
$ws.Range("E20").Value = 16.28000068664551
$ws.Range("F20").Value = 16.43000030517578
$ws.Range("G20").Value = 12.67500019073486
$ws.Range("H20").Value = 328636590
$ws.Range("I20").Value = "PSTG"

$ws.Range("D21").Value = 19.54999923706055
$ws.Range("E21").Value = 16.06999969482422
$ws.Range("F21").Value = 20.5
$ws.Range("G21").Value = 15.80000019073486
$ws.Range("H21").Value = 328636590
$ws.Range("I21").Value = "PSTG"

$ws.Range("D22").Value = 17.84000015258789
$ws.Range("E22").Value = 15.26000022888184
$ws.Range("F22").Value = 19.75
$ws.Range("G22").Value = 14.65999984741211
$ws.Range("H22").Value = 328636590
$ws.Range("I22").Value = "PSTG"

$ws.Range("D23").Value = 13.96000003814697
$ws.Range("E23").Value = 17.61000061035156
$ws.Range("F23").Value = 18.07999992370605
$ws.Range("G23").Value = 12.5649995803833
$ws.Range("H23").Value = 328636590
$ws.Range("I23").Value = "PSTG"

$ws.Range("D24").Value = 18.01000022888184
$ws.Range("E24").Value = 15.26000022888184
$ws.Range("F24").Value = 18.47999954223633
$ws.Range("G24").Value = 14.34000015258789
$ws.Range("H24").Value = 328636590
$ws.Range("I24").Value = "PSTG"

$ws.Range("D25").Value = 16.29000091552734
$ws.Range("E25").Value = 18.27000045776367
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 16.15500068664551
$ws.Range("H25").Value = 328636590
$ws.Range("I25").Value = "PSTG"

$ws.Range("D26").Value = 23.31999969482422
$ws.Range("E26").Value = 23.3799991607666
$ws.Range("F26").Value = 29.53000068664551
$ws.Range("G26").Value = 22.52000045776367
$ws.Range("H26").Value = 328636590
$ws.Range("I26").Value = "PSTG"

$ws.Range("D27").Value = 20.29999923706055
$ws.Range("E27").Value = 19.04999923706055
$ws.Range("F27").Value = 20.3700008392334
$ws.Range("G27").Value = 16.79000091552734
$ws.Range("H27").Value = 328636590
$ws.Range("I27").Value = "PSTG"

$ws.Range("D28").Value = 19.64999961853028
$ws.Range("E28").Value = 25.82999992370605
$ws.Range("F28").Value = 26.09499931335449
$ws.Range("G28").Value = 18.75
$ws.Range("H28").Value = 328636590
$ws.Range("I28").Value = "PSTG"

$ws.Range("D29").Value = 26.85000038146973
$ws.Range("E29").Value = 30.96999931335449
$ws.Range("F29").Value = 32.55799865722656
$ws.Range("G29").Value = 26.53000068664551
$ws.Range("H29").Value = 328636590
$ws.Range("I29").Value = "PSTG"

$ws.Range("D30").Value = 26.71999931335449
$ws.Range("E30").Value = 25.94000053405762
$ws.Range("F30").Value = 28.1200008392334
$ws.Range("G30").Value = 23.34000015258789
$ws.Range("H30").Value = 328636590
$ws.Range("I30").Value = "PSTG"

$ws.Range("D31").Value = 29.29999923706055
$ws.Range("E31").Value = 23.72999954223633
$ws.Range("F31").Value = 30.04000091552734
$ws.Range("G31").Value = 21.89500045776367
$ws.Range("H31").Value = 328636590
$ws.Range("I31").Value = "PSTG"

$ws.Range("D32").Value = 28.30999946594238
$ws.Range("E32").Value = 28.96999931335449
$ws.Range("F32").Value = 31.52000045776367
$ws.Range("G32").Value = 27.77000045776367
$ws.Range("H32").Value = 328636590
$ws.Range("I32").Value = "PSTG"

$ws.Range("D33").Value = 31.39999961853028
$ws.Range("E33").Value = 29.19000053405762
$ws.Range("F33").Value = 32.44800186157227
$ws.Range("G33").Value = 27.32999992370605
$ws.Range("H33").Value = 328636590
$ws.Range("I33").Value = "PSTG"

$ws.Range("D34").Value = 29.02000045776367
$ws.Range("E34").Value = 28.54000091552734
$ws.Range("F34").Value = 32.18000030517578
$ws.Range("G34").Value = 27.96999931335449
$ws.Range("H34").Value = 328636590
$ws.Range("I34").Value = "PSTG"

$ws.Range("D35").Value = 22.90999984741211
$ws.Range("E35").Value = 28.79000091552734
$ws.Range("F35").Value = 29.65999984741211
$ws.Range("G35").Value = 22.13999938964844
$ws.Range("H35").Value = 328636590
$ws.Range("I35").Value = "PSTG"

$ws.Range("D36").Value = 36.58000183105469
$ws.Range("E36").Value = 36.59000015258789
$ws.Range("F36").Value = 38.4900016784668
$ws.Range("G36").Value = 34.0099983215332
$ws.Range("H36").Value = 328636590
$ws.Range("I36").Value = "PSTG"

$ws.Range("D37").Value = 40.36999893188477
$ws.Range("E37").Value = 52.65000152587891
$ws.Range("F37").Value = 52.7599983215332
$ws.Range("G37").Value = 38.77999877929688
$ws.Range("H37").Value = 328636590
$ws.Range("I37").Value = "PSTG"

$ws.Range("D38").Value = 50.33000183105469
$ws.Range("E38").Value = 60.29000091552734
$ws.Range("F38").Value = 68.75
$ws.Range("G38").Value = 49.45000076293945
$ws.Range("H38").Value = 328636590
$ws.Range("I38").Value = "PSTG"

$ws.Range("D39").Value = 59.38000106811523
$ws.Range("E39").Value = 51.29000091552734
$ws.Range("F39").Value = 63.09999847412109
$ws.Range("G39").Value = 49.79000091552734
$ws.Range("H39").Value = 328636590
$ws.Range("I39").Value = "PSTG"

$ws.Range("D40").Value = 50.34000015258789
$ws.Range("E40").Value = 52.9900016784668
$ws.Range("F40").Value = 54.79600143432617
$ws.Range("G40").Value = 45.15000152587891
$ws.Range("H40").Value = 328636590
$ws.Range("I40").Value = "PSTG"

$ws.Range("D41").Value = 65.84999847412109
$ws.Range("E41").Value = 52.47000122070312
$ws.Range("F41").Value = 71.30000305175781
$ws.Range("G41").Value = 50.9010009765625
$ws.Range("H41").Value = 328636590
$ws.Range("I41").Value = "PSTG"

$ws.Range("D42").Value = 46.88000106811523
$ws.Range("E42").Value = 53.59000015258789
$ws.Range("F42").Value = 57.15499877929688
$ws.Range("G42").Value = 46.2599983215332
$ws.Range("H42").Value = 328636590
$ws.Range("I42").Value = "PSTG"

$ws.Range("D43").Value = 57.5
$ws.Range("E43").Value = 77.61000061035156
$ws.Range("F43").Value = 80.68000030517578
$ws.Range("G43").Value = 54.36999893188477
$ws.Range("H43").Value = 328636590
$ws.Range("I43").Value = "PSTG"

